# Update mappings.xlsx with new SupplyLookup entries
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SupplyLookupMappings")
$ws.Activate()

# Row 33 raw value changes from "VL-MG" to "VALERO-MG " (still maps to Valero)
$ws.Cells.Item(33, 1).Value = "VALERO-MG "
$ws.Cells.Item(33, 2).Value = "Valero"
$ws.Rows.Item(33).RowHeight = 13.8

# Insert a new row at 34 so the old row 34 (BP-KM -> BP) shifts down to row 35,
# and put the "VL-MG -> Valero" mapping back in as its own row.
$ws.Rows.Item(34).Insert()
$ws.Cells.Item(34, 1).Value = "VL-MG"
$ws.Cells.Item(34, 2).Value = "Valero"
$ws.Rows.Item(34).RowHeight = 12.8

# New rows 36-37: additional supplier lookup entries
$ws.Cells.Item(36, 1).Value = "KMEP "
$ws.Cells.Item(36, 2).Value = "KMEP "
$ws.Rows.Item(36).RowHeight = 12.8
$ws.Range("A36:B36").WrapText = $true

$ws.Cells.Item(37, 1).Value = "JDS/STL/JDS "
$ws.Cells.Item(37, 2).Value = "JDS Energy "
$ws.Rows.Item(37).RowHeight = 12.8
$ws.Range("A37:B37").WrapText = $true

# Update the sheet's active selection
$ws.Range("G20").Select()
